$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both contain the same event rows in columns
# A-I; column F holds the "想去人数" (want-to-go count) which was refreshed
# for three events (rows 3, 4, 6).
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 2426
    $ws.Range("F4").Value = 454
    $ws.Range("F6").Value = 6514
}
